$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3257249711502525
$ws.Range("H2").Value = 162.7960558837417
$ws.Range("I2").Value = [double]"1.090967547286175e-30"
$ws.Range("J2").Value = 37.84691813850442
$ws.Range("L2").Value = 18.28287537182007
$ws.Range("M2").Value = 0.1123053950697461
$ws.Range("O2").Value = '{''const'': 0.8041965056596425, ''N1ratio-ArgsPreds'': -0.22818582580893026}'
$ws.Range("P2").Value = '{''const'': 1.5299431630824037e-43, ''N1ratio-ArgsPreds'': 1.090967547286159e-30}'
$ws.Range("Q2").Value = '{''N1ratio-ArgsPreds'': -0.5707232001156536}'
$ws.Range("R2").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5707232001156564)}'
$ws.Range("S2").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5707232001156562)}'
$ws.Range("T2").Value = '{''N1ratio-ArgsPreds'': np.float64(32.57249711502554)}'
$ws.Range("G3").Value = 0.4297222340761711
$ws.Range("H3").Value = 62.91987639958763
$ws.Range("I3").Value = [double]"1.343846501666787e-39"
$ws.Range("J3").Value = 32.00957324483368
$ws.Range("L3").Value = 6.030055066372702
$ws.Range("M3").Value = 0.09583704564321462
$ws.Range("O3").Value = '{''const'': 0.6460951361165155, ''N1ratio-ArgsPreds'': -0.2314807109302169, ''latitude'': 0.005203641687904243, ''longitude'': -0.0001992181028474547, ''Macro_class'': 0.05827317632261819}'
$ws.Range("P3").Value = '{''const'': 3.6013124999025194e-29, ''N1ratio-ArgsPreds'': 1.047041916647843e-31, ''latitude'': 1.7345666833023843e-07, ''longitude'': 0.3693533121444341, ''Macro_class'': 3.774054689809112e-09}'
$ws.Range("Q3").Value = '{''N1ratio-ArgsPreds'': -0.5789641474828611, ''latitude'': 0.23003147484609276, ''longitude'': -0.04088923999541228, ''Macro_class'': 0.2551047618687475}'
$ws.Range("R3").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5808785653454406), ''latitude'': np.float64(0.28039353684640295), ''longitude'': np.float64(-0.049126251090201724), ''Macro_class'': np.float64(0.3145097411962766)}'
$ws.Range("S3").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5389014838131498), ''latitude'': np.float64(0.2205931154493497), ''longitude'': np.float64(-0.03714338992437315), ''Macro_class'': np.float64(0.25020430209163114)}'
$ws.Range("T3").Value = '{''N1ratio-ArgsPreds'': np.float64(29.04148092560146), ''latitude'': np.float64(4.866132258365012), ''longitude'': np.float64(0.13796314150740252), ''Macro_class'': np.float64(6.260219278516022)}'
$ws.Range("U3").Value = 0.1039972629259186
$ws.Range("V3").Value = 20.30302184467789
$ws.Range("W3").Value = [double]"4.134710699016487e-12"
$ws.Range("G4").Value = 0.4297976308001815
$ws.Range("H4").Value = 50.20063710268639
$ws.Range("I4").Value = [double]"1.114892836579535e-38"
$ws.Range("J4").Value = 32.00534124228362
$ws.Range("L4").Value = 4.824890453608174
$ws.Range("M4").Value = 0.09611213586271357
$ws.Range("O4").Value = '{''const'': 0.6478715557735011, ''N1ratio-ArgsPreds'': -0.23029823450386092, ''latitude'': 0.005210450677401761, ''longitude'': -0.0002120248347705356, ''Macro_class'': 0.0585405949575478, ''Fam_class'': -0.00010656470313954815}'
$ws.Range("P4").Value = '{''const'': 1.2519041382639054e-28, ''N1ratio-ArgsPreds'': 4.160768616548865e-29, ''latitude'': 1.769814041948919e-07, ''longitude'': 0.3576527504768595, ''Macro_class'': 4.559021367809855e-09, ''Fam_class'': 0.8339225463295897}'
$ws.Range("Q4").Value = '{''N1ratio-ArgsPreds'': -0.5760066161475172, ''latitude'': 0.2303324721072928, ''longitude'': -0.04351780400478214, ''Macro_class'': 0.25627545088019416, ''Fam_class'': -0.009322107485196199}'
$ws.Range("R4").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5605078062590266), ''latitude'': np.float64(0.28060584710727093), ''longitude'': np.float64(-0.0504130760596281), ''Macro_class'': np.float64(0.3133713427156236), ''Fam_class'': np.float64(-0.01149828365769737)}'
$ws.Range("S4").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5110785084455699), ''latitude'': np.float64(0.22075982091380303), ''longitude'': np.float64(-0.03811626027475098), ''Macro_class'': np.float64(0.24918337091065954), ''Fam_class'': np.float64(-0.00868312869941237)}'
$ws.Range("T4").Value = '{''N1ratio-ArgsPreds'': np.float64(26.120124179494848), ''latitude'': np.float64(4.873489852989438), ''longitude'': np.float64(0.14528492973325596), ''Macro_class'': np.float64(6.209235233839933), ''Fam_class'': np.float64(0.007539672401055876)}'
$ws.Range("U4").Value = [double]"7.539672401046804e-05"
$ws.Range("V4").Value = 0.0440319269993904
$ws.Range("W4").Value = 0.833922546329686
$ws.Range("G5").Value = 0.4307960039466122
$ws.Range("H5").Value = 35.78769105502964
$ws.Range("I5").Value = [double]"4.246057926749762e-37"
$ws.Range("J5").Value = 31.94930276372821
$ws.Range("L5").Value = 3.454355820942326
$ws.Range("M5").Value = 0.09652357330431482
$ws.Range("O5").Value = '{''const'': 0.623312746996375, ''N1ratio-ArgsPreds'': -0.22791442762176128, ''latitude'': 0.005240558481256999, ''longitude'': -0.00020654660329224965, ''Macro_class'': 0.055688107563884884, ''Fam_class'': -0.00021841796015278017, ''Nlen_freq'': -0.010433171531082964, ''Vlen_freq'': 0.01555204247010605}'
$ws.Range("P5").Value = '{''const'': 2.9487034745625933e-09, ''N1ratio-ArgsPreds'': 1.0481199564811615e-27, ''latitude'': 1.699262448169016e-07, ''longitude'': 0.3813401191945358, ''Macro_class'': 1.8222686347441684e-07, ''Fam_class'': 0.6813765951816662, ''Nlen_freq'': 0.6739099236055167, ''Vlen_freq'': 0.4791952800585929}'
$ws.Range("Q5").Value = '{''N1ratio-ArgsPreds'': -0.5700443970333957, ''latitude'': 0.23166341357878115, ''longitude'': -0.042393404572881224, ''Macro_class'': 0.2437880053140677, ''Fam_class'': -0.019106849090314564, ''Nlen_freq'': -0.0317491461277499, ''Vlen_freq'': 0.056714222012008236}'
$ws.Range("R5").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.5498644558342528), ''latitude'': np.float64(0.281808191839382), ''longitude'': np.float64(-0.048126418204418056), ''Macro_class'': np.float64(0.2811394426544442), ''Fam_class'': np.float64(-0.02258198396726297), ''Nlen_freq'': np.float64(-0.023142972856167343), ''Vlen_freq'': np.float64(0.038907718982092054)}'
$ws.Range("S5").Value = '{''N1ratio-ArgsPreds'': np.float64(-0.4966736319771695), ''latitude'': np.float64(0.22159289869439136), ''longitude'': np.float64(-0.03635139152273897), ''Macro_class'': np.float64(0.22102185575082223), ''Fam_class'': np.float64(-0.01704146105947908), ''Nlen_freq'': np.float64(-0.01746503461640957), ''Vlen_freq'': np.float64(0.029376409226563177)}'
$ws.Range("T5").Value = '{''N1ratio-ArgsPreds'': np.float64(24.66846967013928), ''latitude'': np.float64(4.910341275178279), ''longitude'': np.float64(0.13214236656394585), ''Macro_class'': np.float64(4.885066071953727), ''Fam_class'': np.float64(0.02904113950417418), ''Nlen_freq'': np.float64(0.03050274341523846), ''Vlen_freq'': np.float64(0.08629734190465062)}'
$ws.Range("U5").Value = 0.0009983731464306977
$ws.Range("V5").Value = 0.2902838997616293
$ws.Range("W5").Value = 0.7482414037482997
